# Rewrite schedule for Autumn 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 38 ---
$ws.Cells.Item(38, 3).Value = "Statistics"
$ws.Cells.Item(38, 4).Value = "Lecture"
$ws.Cells.Item(38, 5).Value = "AI/ChatGPT workshop"

# --- Row 39 ---
$ws.Cells.Item(39, 3).Value = "Statistics"
$ws.Cells.Item(39, 4).Value = "Practical"
$ws.Cells.Item(39, 5).Value = "AI/ChatGPT workshop"

# --- Row 40 ---
$ws.Cells.Item(40, 3).Value = "Assignment"
$ws.Cells.Item(40, 4).Value = "Practical"
$ws.Cells.Item(40, 5).Value = "Intro to Written Assignment"

# --- Row 41 ---
$ws.Cells.Item(41, 3).Value = "Assignment"
$ws.Cells.Item(41, 4).Value = "Practical"
$ws.Cells.Item(41, 5).Value = "Written Assignment 2"

# --- Row 42 ---
$ws.Cells.Item(42, 3).Value = "Assignment"
$ws.Cells.Item(42, 4).Value = "Practical"
$ws.Cells.Item(42, 5).Value = "Written Assignment 3"

# --- Row 43 ---
$ws.Cells.Item(43, 3).Value = "Assignment"
$ws.Cells.Item(43, 4).Value = "Practical"
$ws.Cells.Item(43, 5).Value = "Written Assignment 4"

# --- Row 44 ---
$ws.Cells.Item(44, 3).Value = "Assignment"
$ws.Cells.Item(44, 4).Value = "Practical"
$ws.Cells.Item(44, 5).Value = "Written Assignment 5"

# --- Row 45 ---
$ws.Cells.Item(45, 3).Value = "Assignment"
$ws.Cells.Item(45, 4).Value = "Practical"
$ws.Cells.Item(45, 5).Value = "Written Assignment 6"

# --- New row 46 (copy formatting of row 45 across A:F first) ---
$ws.Range("A45:F45").Copy() | Out-Null
$ws.Range("A46:F46").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(46, 1).Value = 23
$ws.Cells.Item(46, 2).Value = 45
$ws.Cells.Item(46, 3).Value = "Assignment"
$ws.Cells.Item(46, 4).Value = "Practical"
$ws.Cells.Item(46, 5).Value = "Written Assignment 7"
$ws.Cells.Item(46, 6).Value = "OJ"

# --- New row 47 ---
$ws.Range("A45:F45").Copy() | Out-Null
$ws.Range("A47:F47").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(47, 1).Value = 23
$ws.Cells.Item(47, 2).Value = 46
$ws.Cells.Item(47, 3).Value = "Assignment"
$ws.Cells.Item(47, 4).Value = "Practical"
$ws.Cells.Item(47, 5).Value = "Written Assignment 8"
$ws.Cells.Item(47, 6).Value = "OJ"

# --- New row 48 ---
$ws.Range("A45:F45").Copy() | Out-Null
$ws.Range("A48:F48").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(48, 1).Value = 24
$ws.Cells.Item(48, 2).Value = 47
$ws.Cells.Item(48, 3).Value = "Exam"
$ws.Cells.Item(48, 4).Value = "Practical"
$ws.Cells.Item(48, 5).Value = "MCQ Exam"
$ws.Cells.Item(48, 6).Value = "OJ"

# --- New row 49 ---
$ws.Range("A45:F45").Copy() | Out-Null
$ws.Range("A49:F49").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(49, 1).Value = 24
$ws.Cells.Item(49, 2).Value = 48
$ws.Cells.Item(49, 3).Value = "Exam"
$ws.Cells.Item(49, 4).Value = "Practical"
$ws.Cells.Item(49, 5).Value = "MCQ Exam"
$ws.Cells.Item(49, 6).Value = "OJ"

# --- Update sheet view: scroll position & selection ---
$ws.Range("C40:C47").Select() | Out-Null
